$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.589.74"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "3.172.12"
$ws.Range("E3").Value = "  -4.80%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.72"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.40"
$ws.Range("E6").Value = "  -3.17%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("D9").Value = "3.178.55"
$ws.Range("E9").Value = "  -4.50%  "
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.55"
$ws.Range("E11").Value = "  -4.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.391"
$ws.Range("E12").Value = "  -4.46%  "
$ws.Range("D13").Value = "3.727.22"
$ws.Range("E13").Value = "  -4.86%  "
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.26"
$ws.Range("E15").Value = "  -4.87%  "
$ws.Range("D16").Value = "65.555.84"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000163"
$ws.Range("E17").Value = "  -3.41%  "
$ws.Range("D18").Value = "3.181.59"
$ws.Range("E18").Value = "  -4.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.70"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.84"
$ws.Range("E20").Value = "  -4.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "358.34"
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.26"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.05"
$ws.Range("E24").Value = "  -3.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.494"
$ws.Range("E25").Value = "  -4.79%  "
$ws.Range("D26").Value = "3.303.19"
$ws.Range("E26").Value = "  -5.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("E27").Value = "  -5.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.79"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("E29").Value = "  -0.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.92"
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.35"
$ws.Range("E33").Value = "  -5.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.92"
$ws.Range("E34").Value = "  -3.79%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.60"
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.82"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.832"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.42"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.78"
$ws.Range("E41").Value = "  +2.26%  "
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").Value = "2.646.36"
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.10"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.18"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.57"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0657"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.01"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "327.42"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0273"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("E51").Value = "  -1.14%  "
